$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = '44.692.84'
    "D3" = '2.245.74'
    "E3" = '  +0.19%  '
    "E4" = '  +0.08%  '
    "D5" = '''306.44'
    "E5" = '  -0.12%  '
    "D6" = '''94.13'
    "E6" = '  -0.39%  '
    "E7" = '  -0.25%  '
    "E8" = '  +0.02%  '
    "E9" = '  -2.12%  '
    "E10" = '  +0.04%  '
    "D11" = '''0.0799'
    "E11" = '  -1.32%  '
    "D12" = '''7.16'
    "E12" = '  -0.18%  '
    "D13" = '''0.105'
    "E13" = '  +0.14%  '
    "D14" = '2.588.96'
    "E14" = '  +0.24%  '
    "D15" = '2.337.26'
    "E15" = '  -2.16%  '
    "E16" = '  -0.02%  '
    "D17" = '''13.51'
    "E17" = '  +0.09%  '
    "D18" = '44.462.81'
    "E18" = '  +0.97%  '
    "D19" = '0.0₃0934'
    "E19" = '  -3.02%  '
    "D20" = '''6.16'
    "E20" = '  -3.55%  '
    "D21" = '''11.74'
    "E21" = '  -3.18%  '
    "D22" = '''65.29'
    "E22" = '  -0.32%  '
    "D23" = '''237.70'
    "E23" = '  -0.01%  '
    "E24" = '  -0.05%  '
    "E25" = '  -1.46%  '
    "E26" = '  -0.04%  '
    "E27" = '  +3.44%  '
    "E28" = '  -1.49%  '
    "D29" = '''37.03'
    "E29" = '  -4.01%  '
    "D30" = '''19.92'
    "E30" = '  -0.42%  '
    "D31" = '''5.82'
    "E31" = '  -0.16%  '
    "D32" = '''148.62'
    "E32" = '  -2.92%  '
    "E33" = '  +0.37%  '
    "D34" = '''0.0778'
    "E34" = '  -1.94%  '
    "D35" = '''3.18'
    "E35" = '  +1.71%  '
    "D36" = '''0.109'
    "E36" = '  +1.61%  '
    "E37" = '  -2.00%  '
    "E38" = '  +5.12%  '
    "D39" = '''15.03'
    "E39" = '  +5.30%  '
    "E40" = '  -4.55%  '
    "E41" = '  -1.03%  '
    "D42" = '''0.0297'
    "E42" = '  -0.65%  '
    "E43" = '  -0.07%  '
    "D44" = '1.813.50'
    "E44" = '  +3.88%  '
    "E45" = '  +13.52%  '
    "D46" = '''81.48'
    "E46" = '  -1.73%  '
    "E47" = '  -2.43%  '
    "D48" = '''97.80'
    "E48" = '  -1.84%  '
    "E49" = '  -2.17%  '
    "E50" = '  +2.48%  '
    "D51" = '''53.84'
    "E51" = '  -1.48%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
